{"js": "// AWGA_FinalReport.docx \u2014 apply the three \"Project Challenges\" / process\n// section wording edits described by the commit:\n//   1. \"...Clear comments and discussion became key...\"\n//        -> \"...Clear comments and communication became key...\"\n//   2. Rewrite the \"Coordinating with the client...\" paragraph.\n//   3. \"...thought-out the project and not fall too far behind.\"\n//        -> \"...thought-out the project, manage bugs, and not fall too far behind.\"\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) \"discussion\" -> \"communication\" in the collaborative-development\n//    paragraph (scope the search to the surrounding phrase so we don't\n//    touch any other occurrence of the word in the document).\n// ---------------------------------------------------------------------\n{\n  const oldPhrase =\n    \"Clear comments and discussion became key to overcoming this obstacle.\";\n  const newPhrase =\n    \"Clear comments and communication became key to overcoming this obstacle.\";\n\n  const results = body.search(oldPhrase, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newPhrase, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 2) Rework the \"Coordinating with the client\" paragraph.\n// ---------------------------------------------------------------------\n{\n  const oldPara =\n    \"Coordinating with the client was ok at first, but later down the \" +\n    \"line they were not always responsive.  They do after all have real \" +\n    \"jobs.  We were generally able to continue working until they could \" +\n    \"respond.\";\n  const newPara =\n    \"Coordinating with the client went well in the beginning, but later \" +\n    \"on travel and work responsibilities made it difficult to \" +\n    \"communicate effectively.  We were generally able to continue \" +\n    \"working until a response was received.\";\n\n  const results = body.search(oldPara, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newPara, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3) \"thought-out the project and not fall too far behind.\" ->\n//    \"thought-out the project, manage bugs, and not fall too far behind.\"\n//    (the new text lands ahead of the existing _GoBack bookmark, which is\n//    what naturally happens when the replaced range sits before it).\n// ---------------------------------------------------------------------\n{\n  const oldPhrase = \"thought-out the project\";\n  const newPhrase = \"thought-out the project, manage bugs,\";\n\n  const results = body.search(oldPhrase, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newPhrase, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# AWGA_FinalReport.docx \u2014 apply the three \"Project Challenges\" / process\n# section wording edits described by the commit:\n#   1. \"...Clear comments and discussion became key...\"\n#        -> \"...Clear comments and communication became key...\"\n#   2. Rewrite the \"Coordinating with the client...\" paragraph.\n#   3. \"...thought-out the project and not fall too far behind.\"\n#        -> \"...thought-out the project, manage bugs, and not fall too far behind.\"\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"discussion\" -> \"communication\" in the collaborative-development\n#    paragraph (scope the Find to the surrounding phrase so we don't\n#    touch any other occurrence of the word in the document).\n# ---------------------------------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Execute(\n    \"Clear comments and discussion became key to overcoming this obstacle.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Clear comments and communication became key to overcoming this obstacle.\",\n    2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 2) Rework the \"Coordinating with the client\" paragraph.\n# ---------------------------------------------------------------------\n$oldPara = \"Coordinating with the client was ok at first, but later down the line they were not always responsive.  They do after all have real jobs.  We were generally able to continue working until they could respond.\"\n$newPara = \"Coordinating with the client went well in the beginning, but later on travel and work responsibilities made it difficult to communicate effectively.  We were generally able to continue working until a response was received.\"\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Execute($oldPara, $false, $false, $false, $false, $false, $true, 1, $false, $newPara, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3) \"thought-out the project and not fall too far behind.\" ->\n#    \"thought-out the project, manage bugs, and not fall too far behind.\"\n#\n#    \"out the project\" alone is ambiguous (it also matches inside \"tested\n#    throughout the project\" earlier in the document), so first locate the\n#    unique \"thought-\" anchor and scope the real replace to the text that\n#    follows it.\n#\n#    The existing \"_GoBack\" bookmark sits collapsed right in the middle of\n#    \"thought-out the project\" (between \"thought-\" and \"out the project\").\n#    A Find/Replace whose match spans that bookmark drops it, so the text\n#    is rewritten first and the bookmark is then re-added, collapsed, at\n#    its new home: right after the inserted \", manage bugs,\" and before\n#    \" and not fall too far behind\" \u2014 matching where it ends up once the\n#    new text is typed ahead of it.\n# ---------------------------------------------------------------------\n$anchor = $d.Content\n$anchor.Find.ClearFormatting()\n$anchor.Find.Execute(\"thought-\") | Out-Null\n\n$rng3 = $d.Range($anchor.End, $d.Content.End)\n$rng3.Find.ClearFormatting()\n$rng3.Find.Replacement.ClearFormatting()\n$rng3.Find.Execute(\"out the project\", $false, $false, $false, $false, $false, $true, 1, $false, \"out the project, manage bugs,\", 1) | Out-Null\n\n$rng4 = $d.Range($anchor.End, $d.Content.End)\n$rng4.Find.ClearFormatting()\n$rng4.Find.Execute(\"manage bugs,\") | Out-Null\n$bmRange = $rng4.Duplicate\n$bmRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n"}
